# "rebuild to get pages working again"
#
# Renames the two data sheets, bumps the ValueSet metadata (version/date/
# contact) on the Metadata sheet, and inserts a new "Jurisdiction" property
# row right after "Contact" (pushing Description/Purpose/Copyright/Immutable
# down by one row).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Metadata
$ws2 = $wb.Worksheets.Item(2)   # Include ValueSets
$ws3 = $wb.Worksheets.Item(3)   # Exclude from Consent Scope Co

# --- Rename sheets ---
$ws2.Name = "Include ValueSet #0"
$ws3.Name = "Exclude #1"

# --- Metadata sheet value updates ---
$ws1.Range("B3").Value2 = "0.2.2"
$ws1.Range("B8").Value2 = "2024-09-11T16:17:59-05:00"
$ws1.Range("B10").Value2 = "MITRE (https://www.mitre.org)"

# --- Insert a "Jurisdiction" row after "Contact" (row 10), before
#     "Description" (row 11), shifting the remaining property rows down ---

# First, prime the row that will be vacated at the bottom (row 15) with the
# same formatting as the row above it, so nothing loses its style.
$ws1.Range("A14:B14").Copy()
$ws1.Range("A15:B15").PasteSpecial(-4122)

# Shift rows 11-14 down into 12-15 (bottom-up so we don't clobber anything).
for ($r = 14; $r -ge 11; $r--) {
    $dst = $r + 1
    $ws1.Cells.Item($dst, 1).Value2 = $ws1.Cells.Item($r, 1).Value2
    $ws1.Cells.Item($dst, 2).Value2 = $ws1.Cells.Item($r, 2).Value2
}

# Row 11 is now free for the new property.
$ws1.Range("A11").Value2 = "Jurisdiction"
$ws1.Range("B11").Value2 = ""

Write-Host "applied ValueSet metadata rebuild"
